# feat: add 2022-Q4 data
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Update the "总计" (summary) sheet: push the existing per-quarter
#    rows down by one and insert the new 2022-Q4 figures at the top of
#    the data block (row 2), giving the new row-5 index cell the same
#    style as the other "A" index cells.
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

$total.Range("A2").Copy()
$total.Range("A5").PasteSpecial(-4122)
$total.Range("A5").Value = 3

$total.Range("B5").Value = "2022-Q1"
$total.Range("C5").Value = 4
$total.Range("D5").Value = 0.35

$total.Range("B4").Value = "2022-Q2"
$total.Range("C4").Value = 1
$total.Range("D4").Value = 0.02

$total.Range("B3").Value = "2022-Q3"
$total.Range("C3").Value = 3
$total.Range("D3").Value = 0.08

$total.Range("B2").Value = "2022-Q4"
$total.Range("C2").Value = 2
$total.Range("D2").Value = 0.05

# ---------------------------------------------------------------------
# 2) Insert a brand-new "2022-Q4" sheet right before "2022-Q3" (i.e.
#    right after "总计"), carrying the per-fund holdings for the
#    quarter. Duplicating the "2022-Q3" sheet (instead of adding a
#    blank one) keeps the header/index-column styling intact.
# ---------------------------------------------------------------------
$q3 = $wb.Worksheets.Item("2022-Q3")
$q3.Copy($q3)
$wb.Worksheets.Item("2022-Q3 (2)").Name = "2022-Q4"
$q4 = $wb.Worksheets.Item("2022-Q4")

# 2022-Q3 has 3 data rows, 2022-Q4 only has 2 - drop the extra one.
$q4.Rows.Item(4).Delete()

$q4.Range("A2").Value = 0
$q4.Range("B2").Value = "'004194"
$q4.Range("B2").Style = "Normal"
$q4.Range("C2").Value = "招商中证1000指数增强A"
$q4.Range("D2").Value = "'2.57"
$q4.Range("D2").Style = "Normal"
$q4.Range("E2").Value = "'94.27"
$q4.Range("E2").Style = "Normal"
$q4.Range("F2").Value = "'1.02"
$q4.Range("F2").Style = "Normal"
$q4.Range("G2").Value = "'0.0262"
$q4.Range("G2").Style = "Normal"
$q4.Range("H2").Value = 10

$q4.Range("A3").Value = 1
$q4.Range("B3").Value = "'004195"
$q4.Range("B3").Style = "Normal"
$q4.Range("C3").Value = "招商中证1000指数增强C"
$q4.Range("D3").Value = "'2.14"
$q4.Range("D3").Style = "Normal"
$q4.Range("E3").Value = "'94.27"
$q4.Range("E3").Style = "Normal"
$q4.Range("F3").Value = "'1.02"
$q4.Range("F3").Style = "Normal"
$q4.Range("G3").Value = "'0.0218"
$q4.Range("G3").Style = "Normal"
$q4.Range("H3").Value = 10

# Restore the originally active tab (2022-Q1) so the saved workbook's
# tabSelected flag matches the pre-edit file.
$wb.Worksheets.Item("2022-Q1").Activate()
